# Apply cryptos.xlsx data refresh (GitHub Actions update) per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column cells hold numeric-looking text (e.g. "6.42", "0.305") that must
# stay plain text exactly as scraped -- force text format before assigning,
# then restore the default "Normal" style so no stray number-format sticks.
$dCells = @("D2", "D3", "D5", "D6", "D8", "D10", "D13", "D14", "D16", "D18", "D19", "D20", "D26", "D27", "D33", "D35", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "62.854.49"
$ws.Range("E2").Value = "  +1.94%  "

# Row 3
$ws.Range("D3").Value = "3.031.57"
$ws.Range("E3").Value = "  +1.09%  "

# Row 4
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").Value = "593.78"
$ws.Range("E5").Value = "  -0.21%  "

# Row 6
$ws.Range("D6").Value = "153.64"
$ws.Range("E6").Value = "  +6.82%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").Value = "3.025.39"
$ws.Range("E8").Value = "  +0.92%  "

# Row 10
$ws.Range("D10").Value = "6.42"
$ws.Range("E10").Value = "  +7.87%  "

# Row 11
$ws.Range("E11").Value = "  +1.80%  "

# Row 12
$ws.Range("E12").Value = "  +0.76%  "

# Row 13
$ws.Range("D13").Value = "0.0000233"
$ws.Range("E13").Value = "  +2.10%  "

# Row 14
$ws.Range("D14").Value = "35.43"
$ws.Range("E14").Value = "  +3.32%  "

# Row 15
$ws.Range("E15").Value = "  +2.48%  "

# Row 16
$ws.Range("D16").Value = "3.533.22"
$ws.Range("E16").Value = "  +1.28%  "

# Row 17
$ws.Range("E17").Value = "  +1.20%  "

# Row 18
$ws.Range("D18").Value = "62.853.91"
$ws.Range("E18").Value = "  +2.07%  "

# Row 19
$ws.Range("D19").Value = "3.032.01"
$ws.Range("E19").Value = "  +1.12%  "

# Row 20
$ws.Range("D20").Value = "451.14"
$ws.Range("E20").Value = "  -0.71%  "

# Row 21
$ws.Range("E21").Value = "  +1.93%  "

# Row 22
$ws.Range("E22").Value = "  +0.90%  "

# Row 23
$ws.Range("E23").Value = "  +1.44%  "

# Row 24
$ws.Range("E24").Value = "  +0.97%  "

# Row 25
$ws.Range("E25").Value = "  +3.60%  "

# Row 26
$ws.Range("D26").Value = "11.02"
$ws.Range("E26").Value = "  +4.87%  "

# Row 27
$ws.Range("D27").Value = "12.30"
$ws.Range("E27").Value = "  +1.47%  "

# Row 28
$ws.Range("E28").Value = "  +0.00%  "

# Row 29
$ws.Range("E29").Value = "  +6.93%  "

# Row 30
$ws.Range("E30").Value = "  +0.76%  "

# Row 31
$ws.Range("E31").Value = "  +0.07%  "

# Row 32
$ws.Range("E32").Value = "  +5.86%  "

# Row 33
$ws.Range("D33").Value = "27.51"
$ws.Range("E33").Value = "  -0.29%  "

# Row 34
$ws.Range("E34").Value = "  +1.95%  "

# Row 35
$ws.Range("D35").Value = "0.0₃0864"
$ws.Range("E35").Value = "  +5.19%  "

# Row 36
$ws.Range("E36").Value = "  +1.56%  "

# Row 37
$ws.Range("E37").Value = "  +2.77%  "

# Row 38
$ws.Range("E38").Value = "  +10.18%  "

# Row 39
$ws.Range("E39").Value = "  +1.24%  "

# Row 40
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "50.49"
$ws.Range("E40").Value = "  +0.19%  "

# Row 41
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "0.128"
$ws.Range("E41").Value = "  +4.94%  "

# Row 42
$ws.Range("D42").Value = "9.10"
$ws.Range("E42").Value = "  -1.16%  "

# Row 43
$ws.Range("D43").Value = "0.305"
$ws.Range("E43").Value = "  +14.01%  "

# Row 44
$ws.Range("D44").Value = "41.67"
$ws.Range("E44").Value = "  +6.73%  "

# Row 45
$ws.Range("D45").Value = "395.09"
$ws.Range("E45").Value = "  -1.09%  "

# Row 46
$ws.Range("D46").Value = "0.0359"
$ws.Range("E46").Value = "  +1.44%  "

# Row 47
$ws.Range("D47").Value = "2.729.95"
$ws.Range("E47").Value = "  +0.32%  "

# Row 48
$ws.Range("D48").Value = "132.22"
$ws.Range("E48").Value = "  -1.04%  "

# Row 49
$ws.Range("E49").Value = "  +0.02%  "

# Row 50
$ws.Range("E50").Value = "  +3.46%  "

# Row 51
$ws.Range("E51").Value = "  +3.48%  "

# Restore default styling on the D-column cells we text-formatted above
foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}
